$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 57, pushing existing rows 57-78 down to 58-79.
$ws.Rows.Item(57).Insert()

# The row that used to be 57 is now row 58; copy its constant columns into
# the freshly inserted row 57, then overwrite the columns that differ
# (D Fecha, J Volumen, K Precio minimo, L Precio maximo, M Precio promedio
# ponderado, P Precio $/Kg).
$ws.Range("A57").Value = $ws.Range("A58").Value2    # Mercado ID
$ws.Range("B57").Value = $ws.Range("B58").Value2    # Mercado
$ws.Range("C57").Value = $ws.Range("C58").Value2    # Region
$ws.Range("D57").Value = 44463                      # Fecha
$ws.Range("D57").NumberFormat = $ws.Range("D58").NumberFormat
$ws.Range("E57").Value = $ws.Range("E58").Value2    # Codreg
$ws.Range("F57").Value = $ws.Range("F58").Value2    # Categoria ID
$ws.Range("G57").Value = $ws.Range("G58").Value2    # Categoria
$ws.Range("H57").Value = $ws.Range("H58").Value2    # Variedad
$ws.Range("I57").Value = $ws.Range("I58").Value2    # Calidad
$ws.Range("J57").Value = 20                         # Volumen
$ws.Range("K57").Value = 11000                      # Precio minimo
$ws.Range("L57").Value = 12000                      # Precio maximo
$ws.Range("M57").Value = 11500                      # Precio promedio ponderado
$ws.Range("N57").Value = $ws.Range("N58").Value2    # Unidad de comercializacion
$ws.Range("O57").Value = $ws.Range("O58").Value2    # Origen
$ws.Range("P57").Value = 3833                       # Precio $/Kg
$ws.Range("Q57").Value = $ws.Range("Q58").Value2    # Kg o Unidades
$ws.Range("R57").Value = $ws.Range("R58").Value2    # Clasificacion
